# Add "SLS A3" row to the "System Framework" sheet, mirroring the existing
# "SLS" row (row 5) but pointing at the new "SLS A3- Report" report.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("System Framework")

# Copy the formatting (styles, borders, etc.) of the template row (row 5)
# down onto the new row (row 7) before filling in values.
$ws1.Range("A5:F5").Copy() | Out-Null
$ws1.Range("A7:F7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws1.Rows.Item(7).RowHeight = $ws1.Rows.Item(5).RowHeight()

# Fill in the new row's values. "D7" (report name) is written before "A7"
# (computation model) so new shared strings are created in the same order
# as the source edit: "SLS A3- Report" first, then "SLS A3".
$ws1.Range("D7").Value = "SLS A3- Report"
$ws1.Range("A7").Value = "SLS A3"
$ws1.Range("B7").Value = $ws1.Range("B5").Value()
$ws1.Range("C7").Value = $ws1.Range("C5").Value()
$ws1.Range("E7").Value = $ws1.Range("E5").Value()
$ws1.Range("F7").Value = $ws1.Range("F5").Value()

# Reflect that "System Framework" is now the active sheet with the newly
# added row selected (matches the author's final cursor position).
$ws1.Activate()
$ws1.Range("A7").Select() | Out-Null
